# Replace the financial figures in data rows 2-6 with corrected values,
# and remove the now-invalid rows 7-9 (keeping only their label columns A-C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2419
$ws.Range("E2").Value = 131
$ws.Range("F2").Value = 131
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 79
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 19
$ws.Range("K2").Value = 3246
$ws.Range("L2").Value = 1610
$ws.Range("M2").Value = 1636
$ws.Range("N2").Value = 1577
$ws.Range("O2").Value = 59
$ws.Range("P2").Value = 230
$ws.Range("Q2").Value = -90
$ws.Range("R2").Value = -127
$ws.Range("S2").Value = 154
$ws.Range("T2").Value = 230
$ws.Range("U2").Value = -319
$ws.Range("V2").Value = 1274
$ws.Range("W2").Value = 5.42
$ws.Range("X2").Value = 3.26
$ws.Range("Y2").Value = 4.35
$ws.Range("Z2").Value = 2.57
$ws.Range("AA2").Value = 98.43
$ws.Range("AB2").Value = 606.15
$ws.Range("AC2").Value = 134
$ws.Range("AD2").Value = 26.88
$ws.Range("AE2").Value = 2985
$ws.Range("AF2").Value = 1.21
$ws.Range("AG2").Value = 43
$ws.Range("AH2").Value = 1.18
$ws.Range("AI2").Value = 37.22
$ws.Range("AJ2").Value = 53926867

# Row 3
$ws.Range("D3").Value = 2503
$ws.Range("E3").Value = -48
$ws.Range("F3").Value = -48
$ws.Range("G3").Value = -42
$ws.Range("H3").Value = -56
$ws.Range("I3").Value = -43
$ws.Range("J3").Value = -13
$ws.Range("K3").Value = 3109
$ws.Range("L3").Value = 1442
$ws.Range("M3").Value = 1668
$ws.Range("N3").Value = 1620
$ws.Range("O3").Value = 47
$ws.Range("P3").Value = 245
$ws.Range("Q3").Value = 240
$ws.Range("R3").Value = 202
$ws.Range("S3").Value = -128
$ws.Range("T3").Value = 75
$ws.Range("U3").Value = 165
$ws.Range("V3").Value = 1074
$ws.Range("W3").Value = -1.91
$ws.Range("X3").Value = -2.22
$ws.Range("Y3").Value = -2.68
$ws.Range("Z3").Value = -1.75
$ws.Range("AA3").Value = 86.47
$ws.Range("AB3").Value = 576.59
$ws.Range("AC3").Value = -75
$ws.Range("AD3").Value = -94.26
$ws.Range("AE3").Value = 2874
$ws.Range("AF3").Value = 2.46
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 57491208

# Row 4
$ws.Range("D4").Value = 1646
$ws.Range("E4").Value = -292
$ws.Range("F4").Value = -292
$ws.Range("G4").Value = -333
$ws.Range("H4").Value = -332
$ws.Range("I4").Value = -324
$ws.Range("J4").Value = -8
$ws.Range("K4").Value = 2524
$ws.Range("L4").Value = 1191
$ws.Range("M4").Value = 1333
$ws.Range("N4").Value = 1295
$ws.Range("O4").Value = 38
$ws.Range("P4").Value = 245
$ws.Range("Q4").Value = 126
$ws.Range("R4").Value = -180
$ws.Range("S4").Value = -262
$ws.Range("T4").Value = 135
$ws.Range("U4").Value = -8
$ws.Range("V4").Value = 818
$ws.Range("W4").Value = -17.77
$ws.Range("X4").Value = -20.16
$ws.Range("Y4").Value = -22.23
$ws.Range("Z4").Value = -11.78
$ws.Range("AA4").Value = 89.34
$ws.Range("AB4").Value = 444.31
$ws.Range("AC4").Value = -564
$ws.Range("AD4").Value = -4.41
$ws.Range("AE4").Value = 2296
$ws.Range("AF4").Value = 1.08
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 57491208

# Row 5
$ws.Range("D5").Value = 1218
$ws.Range("E5").Value = -97
$ws.Range("F5").Value = -97
$ws.Range("G5").Value = -37
$ws.Range("H5").Value = -67
$ws.Range("I5").Value = -66
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 2413
$ws.Range("L5").Value = 867
$ws.Range("M5").Value = 1546
$ws.Range("N5").Value = 1509
$ws.Range("O5").Value = 37
$ws.Range("P5").Value = 400
$ws.Range("Q5").Value = 59
$ws.Range("R5").Value = -61
$ws.Range("S5").Value = -25
$ws.Range("T5").Value = 48
$ws.Range("U5").Value = 11
$ws.Range("V5").Value = 485
$ws.Range("W5").Value = -7.93
$ws.Range("X5").Value = -5.52
$ws.Range("Y5").Value = -4.73
$ws.Range("Z5").Value = -2.72
$ws.Range("AA5").Value = 56.12
$ws.Range("AB5").Value = 293.99
$ws.Range("AC5").Value = -100
$ws.Range("AD5").Value = -12.64
$ws.Range("AE5").Value = 1913
$ws.Range("AF5").Value = 0.66
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 79983352

# Row 6
$ws.Range("D6").Value = 1084
$ws.Range("E6").Value = -232
$ws.Range("F6").Value = -232
$ws.Range("G6").Value = -296
$ws.Range("H6").Value = -293
$ws.Range("I6").Value = -299
$ws.Range("K6").Value = 2055
$ws.Range("L6").Value = 846
$ws.Range("M6").Value = 1209
$ws.Range("N6").Value = 1209
$ws.Range("P6").Value = 400
$ws.Range("Q6").Value = -71
$ws.Range("R6").Value = 274
$ws.Range("S6").Value = -2
$ws.Range("T6").Value = 15
$ws.Range("U6").Value = -87
$ws.Range("V6").Value = 532
$ws.Range("W6").Value = -21.42
$ws.Range("X6").Value = -27.06
$ws.Range("Y6").Value = -21.98
$ws.Range("Z6").Value = -13.14
$ws.Range("AA6").Value = 69.96
$ws.Range("AB6").Value = 230.68
$ws.Range("AC6").Value = -373
$ws.Range("AD6").Value = -4.51
$ws.Range("AE6").Value = 1532
$ws.Range("AF6").Value = 1.1
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 79983352

# Rows 5-6: clear AG:AH (columns removed in diff)
$ws.Range("AG5:AH6").ClearContents()

# Rows 7-9: clear all data columns D:AI (only A,B,C remain)
$ws.Range("D7:AI9").ClearContents()

